# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.947.96"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.876.97"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'0.7401"
$ws.Range("E5").Value = "  -3.97%  "

$ws.Range("D6").Value = "'242.87"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +0.89%  "

$ws.Range("D9").Value = "'0.07205"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "'24.62"
$ws.Range("E10").Value = "  -3.76%  "

$ws.Range("D11").Value = "'0.08340"
$ws.Range("E11").Value = "  -3.11%  "

$ws.Range("D12").Value = "'0.7518"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").Value = "1.895.76"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "'5.411"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").Value = "'92.57"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").Value = "29.958.59"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("D17").Value = "'6.106"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "'249.13"
$ws.Range("E18").Value = "  +2.01%  "

$ws.Range("D19").Value = "'13.56"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'0.000007860"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").Value = "2.142.89"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").Value = "'8.027"
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'0.1550"
$ws.Range("E25").Value = "  -6.07%  "

$ws.Range("D26").Value = "'9.261"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").Value = "'165.02"
$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").Value = "'2.035"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").Value = "'1.520"
$ws.Range("E30").Value = "  +4.18%  "

$ws.Range("D31").Value = "'4.598"
$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("D32").Value = "'1.537"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "'4.283"
$ws.Range("E33").Value = "  +4.66%  "

$ws.Range("D34").Value = "'0.05324"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").Value = "'0.7489"
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "'0.01966"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'0.4557"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").Value = "1.111.28"
$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("D43").Value = "'6.137"
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("D44").Value = "'72.31"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").Value = "'104.10"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("D48").Value = "'1.854"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").Value = "'7.606"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.496"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.039.18"
$ws.Range("E51").Value = "  -1.08%  "
